{"js": "// Update the division-problem table: replace the text in each of the 25\n// \"problem\" cells (5 data rows x 5 columns, the table also has blank rows\n// in between for students to write their work) with the new values from\n// the commit.\n//\n// Mapping is positional (table row/column index) rather than text-based,\n// since some of the old/new values repeat (e.g. \"94\u00f74=\" is both a target\n// of one change and the source of another), which would make a naive\n// find-and-replace ambiguous.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// logicalRow -> actual table row index (every 4th row holds problems;\n// the rows in between are left blank for answers).\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\n// New text for each of the 5 columns, for each of the 5 data rows, in\n// document order.\nconst newValues = [\n  [\"14\u00f76=\", \"83\u00f73=\", \"62\u00f74=\", \"24\u00f73=\", \"94\u00f74=\"],\n  [\"55\u00f72=\", \"66\u00f77=\", \"89\u00f74=\", \"77\u00f74=\", \"21\u00f72=\"],\n  [\"55\u00f76=\", \"76\u00f72=\", \"13\u00f76=\", \"33\u00f76=\", \"95\u00f72=\"],\n  [\"87\u00f77=\", \"15\u00f78=\", \"97\u00f74=\", \"98\u00f75=\", \"17\u00f79=\"],\n  [\"29\u00f73=\", \"52\u00f79=\", \"25\u00f76=\", \"73\u00f73=\", \"23\u00f75=\"],\n];\n\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const rowIndex = dataRowIndexes[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: replace the text in each of the 25\n# \"problem\" cells (5 data rows x 5 columns; the table also has blank rows\n# in between for students to write their work) with the new values from\n# the commit.\n#\n# Mapping is positional (table row/column index) rather than text-based,\n# since some of the old/new values repeat (e.g. \"94\u00f74=\" is both a target\n# of one change and the source of another), which would make a naive\n# find-and-replace ambiguous.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based table row numbers that hold problems (every 4th row; the rows\n# in between are left blank for answers).\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New text for each of the 5 columns, for each of the 5 data rows, in\n# document order.\n$newValues = @(\n    @(\"14\u00f76=\", \"83\u00f73=\", \"62\u00f74=\", \"24\u00f73=\", \"94\u00f74=\"),\n    @(\"55\u00f72=\", \"66\u00f77=\", \"89\u00f74=\", \"77\u00f74=\", \"21\u00f72=\"),\n    @(\"55\u00f76=\", \"76\u00f72=\", \"13\u00f76=\", \"33\u00f76=\", \"95\u00f72=\"),\n    @(\"87\u00f77=\", \"15\u00f78=\", \"97\u00f74=\", \"98\u00f75=\", \"17\u00f79=\"),\n    @(\"29\u00f73=\", \"52\u00f79=\", \"25\u00f76=\", \"73\u00f73=\", \"23\u00f75=\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowIndex = $dataRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
